# Refresh the legacy GSC export: the first day in the export ("2025-10-05")
# was a blank placeholder row (no Views/Impressions data yet) and has now
# been superseded by a real day's data at the end of the range. Remove that
# leading blank row from the "Chart" sheet; Excel's row delete shifts every
# following row up by one (and drops the now-empty trailing row from the
# used range), which reproduces the refreshed export exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
